$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A6:C6").Copy()
$ws.Range("A7:C7").PasteSpecial(-4122)
$ws.Rows.Item(7).RowHeight = $ws.Rows.Item(6).RowHeight

$ws.Range("A7").Value = 43971
$ws.Range("B7").Value = 2
$ws.Range("C7").Value = "Figure of share of sales: new graphics and variables"

$ws.Range("C8").Select()
